$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare 6 new rows (12-17) by copying the formatting of an existing
# data row so the new cells get the same style (s="2") as other rows.
$ws.Range("A4:B4").Copy()
$ws.Range("A12:B17").PasteSpecial(-4122)

# New rows use a row height of 45, matching similarly sized rows elsewhere
# in the sheet.
$ws.Rows(12).RowHeight = 45
$ws.Rows(13).RowHeight = 45
$ws.Rows(14).RowHeight = 45
$ws.Rows(15).RowHeight = 45
$ws.Rows(16).RowHeight = 45
$ws.Rows(17).RowHeight = 45

# Update column A (Шаги / steps) and column B (Ожидаемые результаты / expected
# results) values for every row - rows 1-11 keep the same set of test cases
# (with rows 3 and 4 swapped and a couple of URLs tidied up) while rows 12-17
# are brand new test cases for the other swapi.co list endpoints.
$ws.Cells.Item(1, 1).Value = "Шаги"
$ws.Cells.Item(1, 2).Value = "Ожидаемые результаты"
$ws.Cells.Item(2, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/                                                Нажать request"
$ws.Cells.Item(2, 2).Value = "Код ответа 200 ОК, json файл с полями people, planets, films, species, vehicles, starships - массивы"
$ws.Cells.Item(3, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/people/3/                                               Нажать request"
$ws.Cells.Item(3, 2).Value = "Код ответа 200 ОК, json файл с полями name, height, mass, hair_color, skin_color, eye_color, birth_year, gender, homeworld,films, species, vehicles и starships - пустой массив, created, edited, url"
$ws.Cells.Item(4, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/people/13/                                               Нажать request"
$ws.Cells.Item(4, 2).Value = "Код ответа 200 ОК, json файл с полями как и в предыдущем тест кейсе, только vehicles и starships уже не пустой массив"
$ws.Cells.Item(5, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/films/5/                                               Нажать request"
$ws.Cells.Item(5, 2).Value = "Код ответа 200 ОК, json файл с полями title, episode_id, opening_crawl, director, producer, release_date,массивы: characters, planets, starships, vehicles;  created, edited, url"
$ws.Cells.Item(6, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/films/10/                                               Нажать request"
$ws.Cells.Item(6, 2).Value = "404 error"
$ws.Cells.Item(7, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/planets/15/                                               Нажать request"
$ws.Cells.Item(7, 2).Value = "Код ответа 200 ОК, json файл с полями name, rotation_period, orbital_period, diameter, climate, gravity, terrain, surface_water, population,  массив:  residents, films; created, edited, url"
$ws.Cells.Item(8, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/species/7/                                               Нажать request"
$ws.Cells.Item(8, 2).Value = "Код ответа 200 ОК, json файл с полями name, classification, designation, average_height, skin_colors, hair_colors, eye_colors, average_lifespan, homeworld, language,  массив:  people, films; created, edited, url"
$ws.Cells.Item(9, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/starships/2/                                               Нажать request"
$ws.Cells.Item(9, 2).Value = "Код ответа 200 ОК, json файл с полями name, model, manufacturer, cost_in_credits, length, max_atmosphering_speed, crew, passengers, cargo_capacity, consumables, hyperdrive_rating, MGLT, starship_class  массив:  pilots, films; created, edited, url"
$ws.Cells.Item(10, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/vehicles/4/                                               Нажать request"
$ws.Cells.Item(10, 2).Value = "Код ответа 200 ОК, json файл с полями name, model, manufacturer, cost_in_credits, length, max_atmosphering_speed, crew, passengers, cargo_capacity, consumables, vehicle_class,  массив:  pilots, films; created, edited, url"
$ws.Cells.Item(11, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/name/4/                                               Нажать request"
$ws.Cells.Item(11, 2).Value = "404 error"
$ws.Cells.Item(12, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/people/                                               Нажать request"
$ws.Cells.Item(12, 2).Value = "Код ответа 200 ОК, json файл с полями count, next, previous, массив results с информацией о каждом из people; created, edited, url"
$ws.Cells.Item(13, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/planets/                                               Нажать request"
$ws.Cells.Item(13, 2).Value = "Код ответа 200 ОК, json файл с полями count, next, previous, массив results с информацией о каждой из planets; created, edited, url"
$ws.Cells.Item(14, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/vehicles/                                               Нажать request"
$ws.Cells.Item(14, 2).Value = "Код ответа 200 ОК, json файл с полями count, next, previous, массив results с информацией о каждом из vehicles; created, edited, url"
$ws.Cells.Item(15, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/starships/                                               Нажать request"
$ws.Cells.Item(15, 2).Value = "Код ответа 200 ОК, json файл с полями count, next, previous, массив results с информацией о каждом из starships; created, edited, url"
$ws.Cells.Item(16, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/films/                                               Нажать request"
$ws.Cells.Item(16, 2).Value = "Код ответа 200 ОК, json файл с полями count, next, previous, массив results с информацией о каждом из films; created, edited, url"
$ws.Cells.Item(17, 1).Value = "Зайти на на сайт                                                           Ввести https://swapi.co/api/species/                                               Нажать request"
$ws.Cells.Item(17, 2).Value = "Код ответа 200 ОК, json файл с полями count, next, previous, массив results с информацией о каждом из species; created, edited, url"

# Move the sheet view / selection so the newly added last row is visible,
# matching the author's saved cursor position.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("B17").Select() | Out-Null
